# Báo cáo bán hàng: đổi placeholder từ "IndirectSalesOrders" sang "SalesOrders"
# (thêm báo cáo bán hàng theo đơn trực tiếp)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$fields = @("STT","Code","BuyerStoreName","SellerStoreName","SaleEmployeeName","eOrderDate","SubTotal","Discount","TaxValue","Total")
$cols = @("A","B","C","D","E","F","G","H","I","J")

for ($i = 0; $i -lt $fields.Length; $i++) {
    $addr = "$($cols[$i])9"
    $ws.Range($addr).Value = "{{ReportSalesOrderGenerals.SalesOrders.$($fields[$i])}}"
}
